# Auto-generated Excel COM-interop script
# Applies currentAveragePrice / Leve price / profit updates pulled by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 8266.75
$ws.Range("I47").Value = 6022.3335
$ws.Range("K47").Value = 6022.3335
$ws.Range("M47").Value = -5050.3335
$ws.Range("H64").Value = 3310
$ws.Range("J64").Value = 3310
$ws.Range("L64").Value = 3310
$ws.Range("N64").Value = -3806
$ws.Range("H67").Value = 3310
$ws.Range("J67").Value = 3310
$ws.Range("L67").Value = 3310
$ws.Range("N67").Value = -5026
$ws.Range("H99").Value = 982.44446
$ws.Range("I99").Value = 334.57144
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 1003.71432
$ws.Range("L99").Value = 9750
$ws.Range("M99").Value = 494.28568
$ws.Range("N99").Value = -12746
$ws.Range("H107").Value = 390.57144
$ws.Range("I107").Value = 315.33334
$ws.Range("K107").Value = 315.33334
$ws.Range("M107").Value = 1604.66666
$ws.Range("H132").Value = 1025.6444
$ws.Range("I132").Value = 949.381
$ws.Range("J132").Value = 2093.3333
$ws.Range("K132").Value = 2848.143
$ws.Range("L132").Value = 6279.999899999999
$ws.Range("M132").Value = -318.143
$ws.Range("N132").Value = -11339.9999
$ws.Range("H133").Value = 87133.164
$ws.Range("J133").Value = 87133.164
$ws.Range("L133").Value = 87133.164
$ws.Range("N133").Value = -97253.164
$ws.Range("H135").Value = 641.25
$ws.Range("I135").Value = 590
$ws.Range("K135").Value = 5310
$ws.Range("M135").Value = -2775
$ws.Range("H137").Value = 1773.4166
$ws.Range("I137").Value = 1570.2858
$ws.Range("K137").Value = 4710.857400000001
$ws.Range("M137").Value = -2160.857400000001
$ws.Range("H138").Value = 2692.4443
$ws.Range("I138").Value = 3288.5
$ws.Range("J138").Value = 2341.8235
$ws.Range("K138").Value = 9865.5
$ws.Range("L138").Value = 7025.470499999999
$ws.Range("M138").Value = -4725.5
$ws.Range("N138").Value = -17305.4705
$ws.Range("H140").Value = 49325.855
$ws.Range("J140").Value = 49325.855
$ws.Range("L140").Value = 49325.855
$ws.Range("N140").Value = -59685.855
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2787.2173
$ws.Range("I32").Value = 2050.93
$ws.Range("K32").Value = 2050.93
$ws.Range("M32").Value = -1763.93
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 2225
$ws.Range("I122").Value = 929.1667
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 2787.5001
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -337.5001000000002
$ws.Range("N122").Value = -34900
$ws.Range("H132").Value = 1353.8928
$ws.Range("I132").Value = 1005.3488
$ws.Range("J132").Value = 2506.7693
$ws.Range("K132").Value = 3016.0464
$ws.Range("L132").Value = 7520.3079
$ws.Range("M132").Value = -486.0464000000002
$ws.Range("N132").Value = -12580.3079
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 70000
$ws.Range("J76").Value = 70000
$ws.Range("L76").Value = 70000
$ws.Range("N76").Value = -70630
$ws.Range("H79").Value = 70000
$ws.Range("J79").Value = 70000
$ws.Range("L79").Value = 70000
$ws.Range("N79").Value = -72184
$ws.Range("H107").Value = 3802.6
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 1838.5278
$ws.Range("I134").Value = 1864.9032
$ws.Range("J134").Value = 1675
$ws.Range("K134").Value = 5594.7096
$ws.Range("L134").Value = 5025
$ws.Range("M134").Value = -3059.7096
$ws.Range("N134").Value = -10095
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 70000
$ws.Range("J28").Value = 70000
$ws.Range("L28").Value = 70000
$ws.Range("N28").Value = -70490
$ws.Range("H31").Value = 1843.125
$ws.Range("I31").Value = 1647.7
$ws.Range("J31").Value = 2168.8333
$ws.Range("K31").Value = 1647.7
$ws.Range("L31").Value = 2168.8333
$ws.Range("M31").Value = -1352.7
$ws.Range("N31").Value = -2758.8333
$ws.Range("H34").Value = 1843.125
$ws.Range("I34").Value = 1647.7
$ws.Range("J34").Value = 2168.8333
$ws.Range("K34").Value = 1647.7
$ws.Range("L34").Value = 2168.8333
$ws.Range("M34").Value = -1445.7
$ws.Range("N34").Value = -2572.8333
$ws.Range("H107").Value = 515.1053000000001
$ws.Range("I107").Value = 455.13333
$ws.Range("J107").Value = 740
$ws.Range("K107").Value = 455.13333
$ws.Range("L107").Value = 740
$ws.Range("M107").Value = 1464.86667
$ws.Range("N107").Value = -4580
$ws.Range("H132").Value = 1918.2858
$ws.Range("I132").Value = 1321.6875
$ws.Range("K132").Value = 3965.0625
$ws.Range("M132").Value = -1435.0625
$ws.Range("H134").Value = 1888.32
$ws.Range("I134").Value = 1643.2858
$ws.Range("J134").Value = 3174.75
$ws.Range("K134").Value = 4929.857400000001
$ws.Range("L134").Value = 9524.25
$ws.Range("M134").Value = -2394.857400000001
$ws.Range("N134").Value = -14594.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.08333
$ws.Range("I2").Value = 155.42857
$ws.Range("J2").Value = 142.6
$ws.Range("K2").Value = 932.57142
$ws.Range("L2").Value = 855.5999999999999
$ws.Range("M2").Value = -819.57142
$ws.Range("N2").Value = -1081.6
$ws.Range("H33").Value = 84.666664
$ws.Range("I33").Value = 77.2
$ws.Range("K33").Value = 463.2
$ws.Range("M33").Value = -180.2
$ws.Range("H44").Value = 299
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 299
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 897
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -1693
$ws.Range("H103").Value = 2143.1428
$ws.Range("I103").Value = 2085.8
$ws.Range("J103").Value = 2175
$ws.Range("K103").Value = 6257.400000000001
$ws.Range("L103").Value = 6525
$ws.Range("M103").Value = -5378.400000000001
$ws.Range("N103").Value = -8283
$ws.Range("H107").Value = 797.9474
$ws.Range("J107").Value = 918.3333
$ws.Range("L107").Value = 2754.9999
$ws.Range("N107").Value = -6594.9999
$ws.Range("H119").Value = 999
$ws.Range("I119").Value = 999
$ws.Range("K119").Value = 2997
$ws.Range("M119").Value = 1841
$ws.Range("H131").Value = 10905.531
$ws.Range("J131").Value = 11908.069
$ws.Range("L131").Value = 35724.20699999999
$ws.Range("N131").Value = -45804.20699999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5529666.5
$ws.Range("I7").Value = 5730769
$ws.Range("J7").Value = 5006800
$ws.Range("K7").Value = 5730769
$ws.Range("L7").Value = 5006800
$ws.Range("M7").Value = -5730657
$ws.Range("N7").Value = -5007024
$ws.Range("H8").Value = 5529666.5
$ws.Range("I8").Value = 5730769
$ws.Range("J8").Value = 5006800
$ws.Range("K8").Value = 5730769
$ws.Range("L8").Value = 5006800
$ws.Range("M8").Value = -5730630
$ws.Range("N8").Value = -5007078
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 10000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -19984
$ws.Range("H122").Value = 1807.762
$ws.Range("I122").Value = 1337.6
$ws.Range("K122").Value = 4012.8
$ws.Range("M122").Value = -1562.8
$ws.Range("H132").Value = 1758.1207
$ws.Range("I132").Value = 1409.2094
$ws.Range("K132").Value = 4227.6282
$ws.Range("M132").Value = -1697.6282
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1879.05
$ws.Range("I132").Value = 1621.1765
$ws.Range("K132").Value = 4863.529500000001
$ws.Range("M132").Value = -2333.529500000001
$ws.Range("H136").Value = 3274.5264
$ws.Range("I136").Value = 2530.7778
$ws.Range("J136").Value = 3943.9
$ws.Range("K136").Value = 7592.3334
$ws.Range("L136").Value = 11831.7
$ws.Range("M136").Value = -5042.3334
$ws.Range("N136").Value = -16931.7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H69").Value = 11592
$ws.Range("J69").Value = 11592
$ws.Range("L69").Value = 11592
$ws.Range("N69").Value = -13090
$ws.Range("H72").Value = 11592
$ws.Range("J72").Value = 11592
$ws.Range("L72").Value = 34776
$ws.Range("N72").Value = -42264
$ws.Range("H104").Value = 14444
$ws.Range("J104").Value = 14444
$ws.Range("L104").Value = 14444
$ws.Range("N104").Value = -21432
$ws.Range("H109").Value = 79999
$ws.Range("J109").Value = 79999
$ws.Range("L109").Value = 79999
$ws.Range("N109").Value = -82773
$ws.Range("H126").Value = 3944.6924
$ws.Range("I126").Value = 2128.1
$ws.Range("K126").Value = 6384.299999999999
$ws.Range("M126").Value = -3914.299999999999

Write-Host "Applied cell updates successfully"
